$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column S (2022) is added as a copy of column R's formatting, then values
# are adjusted per row to match the new 2022 data.

# Row 3: header year 2022
$ws.Range("R3").Copy($ws.Range("S3"))
$ws.Range("S3").Value = 2022

# Row 4: numeric 0
$ws.Range("R4").Copy($ws.Range("S4"))
$ws.Range("S4").Value = 0

# Rows 5-10: no data ("-")
$ws.Range("R5").Copy($ws.Range("S5"))
$ws.Range("R6").Copy($ws.Range("S6"))
$ws.Range("R7").Copy($ws.Range("S7"))
$ws.Range("R8").Copy($ws.Range("S8"))
$ws.Range("R9").Copy($ws.Range("S9"))
$ws.Range("R10").Copy($ws.Range("S10"))

# Row 11: numeric 0 (not "-")
$ws.Range("R11").Copy($ws.Range("S11"))
$ws.Range("S11").Value = 0

# Row 12: same value as column R (1.8411781330637848E-3)
$ws.Range("R12").Copy($ws.Range("S12"))

# Row 13: no data ("-")
$ws.Range("R13").Copy($ws.Range("S13"))

# Leave selection on T3, matching the cursor position after the edit.
[void]$ws.Range("T3").Select()
